$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G3").Value = "不可售"
$ws1.Range("G4").Value = 70
$ws1.Range("F6").Value = 522
$ws1.Range("F10").Value = 40
$ws1.Range("F11").Value = 6996
$ws1.Range("F12").Value = 251
$ws1.Range("F13").Value = 398
$ws1.Range("F14").Value = 3430
$ws1.Range("F15").Value = 246
$ws1.Range("F16").Value = 440
$ws1.Range("F18").Value = 580
$ws1.Range("F19").Value = 56

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12
$ws2.Range("F3").Value = 12
$ws2.Range("G3").Value = 178

# Sheet "全部类型" (sheet4) - superset; mirrors the same 演出 rows (2,3)
# plus the same 展览 rows (5,6,8,12,14,16,17,18,19,20,22,23)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12
$ws4.Range("F3").Value = 12
$ws4.Range("G3").Value = 178
$ws4.Range("G5").Value = "不可售"
$ws4.Range("G6").Value = 70
$ws4.Range("F8").Value = 522
$ws4.Range("F12").Value = 40
$ws4.Range("F14").Value = 6996
$ws4.Range("F16").Value = 251
$ws4.Range("F17").Value = 398
$ws4.Range("F18").Value = 3430
$ws4.Range("F19").Value = 246
$ws4.Range("F20").Value = 440
$ws4.Range("F22").Value = 580
$ws4.Range("F23").Value = 56

# Sheet "本地生活" (sheet3) - unchanged by this edit
